# Update countries & provincias Spain
# Applies the COVID data refresh described by the commit:
#  - re-sorted a few neighbouring countries (their case numbers changed
#    rank order), which manifests as the country label AND its row of
#    numbers moving to a different row
#  - updated the "datos actualizados" timestamp string
#  - refreshed a batch of per-country case/recovered/death counters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-sorted country pairs/triples: the country name shown in a given
#    row changes (while the row's position stays put), and that row's
#    statistics get the refreshed numbers for the new country.
# ---------------------------------------------------------------------
$ws.Range("A63").Value = "Moldavia"
$ws.Range("A64").Value = "Uzbekistan"

$ws.Range("A91").Value = "Tayikistan"
$ws.Range("A92").Value = "Haiti"
$ws.Range("A93").Value = "Finlandia"

$ws.Range("A98").Value = "Albania"
$ws.Range("A99").Value = "Paraguay"

# ---------------------------------------------------------------------
# 2) Timestamp footer string
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 17:12"

# ---------------------------------------------------------------------
# 3) Refreshed numeric counters (Casos totales / Nuevos casos /
#    Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes)
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 4774590
$ws.Range("C4").Value = 10272
$ws.Range("D4").Value = 2363229
$ws.Range("E4").Value = 2253364
$ws.Range("G4").Value = 99
$ws.Range("H4").Value = 157997

$ws.Range("B6").Value = 1780268
$ws.Range("C6").Value = 28349
$ws.Range("D6").Value = 1165442
$ws.Range("E6").Value = 577136
$ws.Range("G6").Value = 287
$ws.Range("H6").Value = 37690

$ws.Range("D22").Value = 89026
$ws.Range("E22").Value = 103905
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 3612

$ws.Range("D46").Value = 46926
$ws.Range("E46").Value = 5872

$ws.Range("B63").Value = 25362
$ws.Range("C63").Value = 249
$ws.Range("D63").Value = 17816
$ws.Range("E63").Value = 6755
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 791

$ws.Range("B64").Value = 25336
$ws.Range("C64").Value = 553
$ws.Range("D64").Value = 15833
$ws.Range("E64").Value = 9352
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 151

$ws.Range("B91").Value = 7495
$ws.Range("D91").Value = 6276
$ws.Range("E91").Value = 1158
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 61

$ws.Range("B92").Value = 7468
$ws.Range("C92").Value = 44
$ws.Range("D92").Value = 4606
$ws.Range("E92").Value = 2697
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 165

$ws.Range("B93").Value = 7453
$ws.Range("C93").Value = 10
$ws.Range("D93").Value = 6950
$ws.Range("E93").Value = 174
$ws.Range("H93").Value = 329

$ws.Range("B98").Value = 5519
$ws.Range("C98").Value = 123
$ws.Range("D98").Value = 3018
$ws.Range("E98").Value = 2335
$ws.Range("G98").Value = 5
$ws.Range("H98").Value = 166

$ws.Range("B99").Value = 5485
$ws.Range("D99").Value = 3786
$ws.Range("E99").Value = 1647
$ws.Range("H99").Value = 52

$ws.Range("B120").Value = 2646
$ws.Range("C120").Value = 13
$ws.Range("D120").Value = 2369
$ws.Range("E120").Value = 190

$ws.Range("B180").Value = 180
$ws.Range("C180").Value = 7
$ws.Range("E180").Value = 40
